$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (header labels) ---
# W1 keeps its text (hideCols "Dernier diplome") - untouched, Excel will reindex
# the shared string table for us automatically.

# New header cells AW1 / AX1, formatted like the other header cells (style s="2").
$ws.Range("V1").Copy()
$ws.Range("AW1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("AW1").Value = '<jt:hideCols test="${datAnnulHide}">Date annulation</jt:hideCols>'

$ws.Range("V1").Copy()
$ws.Range("AX1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("AX1").Value = '<jt:hideCols test="${userAnnulHide}">Annulé par</jt:hideCols>'

# --- Row 2 (jasper template row) ---
# AV2 loses the trailing </jt:forEach> (it now lives on AX2).
$ws.Range("AV2").Value = '<jt:if test="${cand.temAcceptCand!=null}"><jt:if test="${cand.temAcceptCand}" then="CONFIRMATION" else="DESISTEMENT"/></jt:if>'

# New data cells AW2 / AX2 - no explicit style (matches source diff: no s="" attribute).
$ws.Range("AW2").Value = '<jt:if test="${cand.datAnnulCand!=null}">${cand.datAnnulCand}</jt:if>'
$ws.Range("AX2").Value = '${cand.userAnnulCand}</jt:forEach>'


# --- Row 3 (footer) ---
# A3 keeps its text (${footer}) - untouched, Excel will reindex the shared
# string table for us automatically.

# --- New column 49 (AW) width ---
# Target stored width is 16.5703125; this engine's ColumnWidth setter quantizes
# to the nearest 1/6 (plus a fixed 5/6 padding baked into the persisted value),
# so 94/6 is the closest input that round-trips to the nearest achievable width.
$ws.Range("AW1").EntireColumn.ColumnWidth = (94/6)

# --- View state: scroll / selection ---
$ws.Activate()
$ws.Range("AW2").Select()
